$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add header cells I1 ("I0") and J1 ("IF") ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold font, border, centered alignment) from the
# existing header cell H1 onto the two new header cells so the new
# headers match the look of the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fill in the I0/IF data for rows 2 through 84 ---
$data = @(
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(7, 8),
    @(8, 8),
    @(8, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(7, 8),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(6, 6),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(7, 8),
    @(9, 9),
    @(7, 8),
    @(8, 8),
    @(9, 10),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(6, 7),
    @(11, 12),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(6, 7),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(6, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(6, 7),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(6, 6),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(4, 4),
    @(4, 4)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $r = $idx + 2
    $pair = $data[$idx]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}
